{"js": "const body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Locate the three consecutive paragraphs that need to be removed:\n//   1) the blank paragraph right after \"LOQ4057: Opera\u00e7\u00f5es Unit\u00e1rias III (Requisito fraco)\"\n//   2) \"Ver no Jupiter Salvar em pdf Salvar em docx\"\n//   3) \"\u00a9 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github pages. Original theme under Creative Commons Attribution\"\n\nconst items = paragraphs.items;\nlet anchorIndex = -1;\nfor (let i = 0; i < items.length; i++) {\n  if (items[i].text.indexOf(\"LOQ4057\") !== -1) {\n    anchorIndex = i;\n    break;\n  }\n}\n\nif (anchorIndex === -1) {\n  throw new Error(\"Could not find the 'LOQ4057' paragraph anchor.\");\n}\n\nconst blankIndex = anchorIndex + 1;\nconst jupiterIndex = anchorIndex + 2;\nconst copyrightIndex = anchorIndex + 3;\n\nif (\n  items[jupiterIndex].text.indexOf(\"Ver no Jupiter\") === -1 ||\n  items[copyrightIndex].text.indexOf(\"Contact:\") === -1\n) {\n  throw new Error(\"Unexpected document structure around the 'LOQ4057' paragraph.\");\n}\n\n// Delete from the bottom up so earlier indices stay valid.\nitems[copyrightIndex].delete();\nitems[jupiterIndex].delete();\nitems[blankIndex].delete();\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# Locate the \"LOQ4057: Opera\u00e7\u00f5es Unit\u00e1rias III (Requisito fraco)\" paragraph;\n# the three paragraphs right after it (a blank paragraph, the \"Ver no\n# Jupiter...\" line, and the \"\u00a9 2020 ...\" footer line) must be removed.\n$count = $d.Paragraphs.Count\n$anchorIndex = -1\nfor ($i = 1; $i -le $count; $i++) {\n    if ($d.Paragraphs.Item($i).Range.Text -like \"*LOQ4057*\") {\n        $anchorIndex = $i\n        break\n    }\n}\n\nif ($anchorIndex -eq -1) {\n    throw \"Could not find the 'LOQ4057' paragraph anchor.\"\n}\n\n$blankIndex = $anchorIndex + 1\n$jupiterIndex = $anchorIndex + 2\n$copyrightIndex = $anchorIndex + 3\n\nif (($d.Paragraphs.Item($jupiterIndex).Range.Text -notlike \"*Ver no Jupiter*\") -or\n    ($d.Paragraphs.Item($copyrightIndex).Range.Text -notlike \"*Contact:*\")) {\n    throw \"Unexpected document structure around the 'LOQ4057' paragraph.\"\n}\n\n# Delete from the bottom up so earlier indices stay valid.\n$d.Paragraphs.Item($copyrightIndex).Range.Delete()\n$d.Paragraphs.Item($jupiterIndex).Range.Delete()\n$d.Paragraphs.Item($blankIndex).Range.Delete()\n"}
